$wb = $excel.ActiveWorkbook

# --- CAOUser: swap out departed team member for new hire ---
$caoUser = $wb.Worksheets.Item("CAOUser")
$caoUser.Range("A2").Value = "Blaise Brunda"
$caoUser.Range("A2").Select()

# --- NewDealTeamMembers: rename deal-team contact + add Counterparty column ---
$newDeal = $wb.Worksheets.Item("NewDealTeamMembers")
$newDeal.Range("A3").Value = "Amy Xia"
$newDeal.Range("G3").Value = "Ashley Choi"
$newDeal.Range("D11").Select()
